# Applies the "final version for testing" edit to the Configuration sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# REPETITIONS: 79 -> 29
$ws.Range("B6").Value = 29

# FRIEND_RECOMMENDATION: 0 -> 1
$ws.Range("B15").Value = 1

# New row 16 (SCENARIO, value 1) - match formatting of the preceding
# label/value pair (A15:B15).
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# New row 17 (LEARNING_PERIODS, value 100) - the label cell A17 keeps the
# default (no explicit) style, while B17 matches the other value cells.
$ws.Range("B15").Copy()
$ws.Range("B17").PasteSpecial(-4122)

$ws.Range("A16").Value = "SCENARIO"
$ws.Range("B16").Value = 1
$ws.Range("A17").Value = "LEARNING_PERIODS"
$ws.Range("B17").Value = 100

# Update the saved selection to match the authored state (B18).
$ws.Range("B18").Select()
